# Weekly refresh of the "Terminal Hortofrutícola Agro Chillán - Alcachofa" sheet.
# Updates several existing data rows (2-14) with refreshed date/volume/price
# figures and appends a brand-new record as row 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 44455
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 13500
$ws.Range("P2").Value = 338

# --- Row 3 ---
$ws.Range("D3").Value = 44420

# --- Row 4 ---
$ws.Range("D4").Value = 44446
$ws.Range("K4").Value = 12500
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 12750
$ws.Range("P4").Value = 319

# --- Row 5 ---
$ws.Range("D5").Value = 44417
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 15500
$ws.Range("P5").Value = 388

# --- Row 6 ---
$ws.Range("D6").Value = 44399
$ws.Range("H6").Value = "Española"
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 15500
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15750
$ws.Range("P6").Value = 394

# --- Row 7 ---
$ws.Range("D7").Value = 44454
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 13500
$ws.Range("P7").Value = 338

# --- Row 9 ---
$ws.Range("D9").Value = 44427
$ws.Range("H9").Value = "Madrigal"
$ws.Range("I9").Value = "Primera"
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 13500
$ws.Range("P9").Value = 338

# --- Row 10 ---
$ws.Range("D10").Value = 44473
$ws.Range("J10").Value = 160
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11500
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 288

# --- Row 11 ---
$ws.Range("D11").Value = 44426
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 13500
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 338

# --- Row 12 ---
$ws.Range("D12").Value = 44475
$ws.Range("K12").Value = 11000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 11500
$ws.Range("P12").Value = 288

# --- Row 13 ---
$ws.Range("D13").Value = 44453
$ws.Range("J13").Value = 160
$ws.Range("K13").Value = 12500
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 12750
$ws.Range("P13").Value = 319

# --- Row 14 ---
$ws.Range("D14").Value = 44435
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14500
$ws.Range("P14").Value = 362

# --- New row 17 (appended record) ---
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 44432
$ws.Range("D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 100112013
$ws.Range("G17").Value = "Alcachofa"
$ws.Range("H17").Value = "Madrigal"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14500
$ws.Range("N17").Value = "$/caja 40 unidades"
$ws.Range("O17").Value = "Provincia del Elquí"
$ws.Range("P17").Value = 362
$ws.Range("Q17").Value = 40
$ws.Range("R17").Value = "Hortaliza"
